$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12
$ws.Range("A12").Value = 131202678
$ws.Range("B12").Value = 57881
$ws.Range("D12").Value = "NT"
$ws.Range("E12").Value = 100049
$ws.Range("F12").Value = "Spillkråka"
$ws.Range("G12").Value = "Dryocopus martius"
$ws.Range("H12").Value = "(Linnaeus, 1758)"
$ws.Range("M12").Value = "äldre spår"
$ws.Range("Q12").Value = 485445
$ws.Range("R12").Value = 6783152

# Row 13
$ws.Range("A13").Value = 131202690
$ws.Range("B13").Value = 8451
$ws.Range("D13").Value = "LC"
$ws.Range("E13").Value = 106545
$ws.Range("F13").Value = "Mindre märgborre"
$ws.Range("G13").Value = "Tomicus minor"
$ws.Range("H13").Value = "(Hartig, 1834)"
$ws.Range("M13").Value = "äldre gnagspår"
$ws.Range("Q13").Value = 485448
$ws.Range("R13").Value = 6783153

# Row 25
$ws.Range("A25").Value = 131198972
$ws.Range("Q25").Value = 485466
$ws.Range("R25").Value = 6783153

# Row 26
$ws.Range("A26").Value = 131202538
$ws.Range("Q26").Value = 485500
$ws.Range("R26").Value = 6783141

# Row 31
$ws.Range("A31").Value = 131198849
$ws.Range("Q31").Value = 485436
$ws.Range("R31").Value = 6783096
$ws.Range("AC31").Value = "rikligt"

# Row 32
$ws.Range("A32").Value = 131202595
$ws.Range("B32").Value = 57884
$ws.Range("D32").Value = "NT"
$ws.Range("E32").Value = 100109
$ws.Range("F32").Value = "Tretåig hackspett"
$ws.Range("G32").Value = "Picoides tridactylus"
$ws.Range("H32").Value = "(Linnaeus, 1758)"
$ws.Range("M32").Value = "färska spår"
$ws.Range("Q32").Value = 485519
$ws.Range("R32").Value = 6783172
$ws.Range("S32").Value = 25

# Row 33
$ws.Range("A33").Value = 131202509
$ws.Range("B33").Value = 8451
$ws.Range("D33").Value = "LC"
$ws.Range("E33").Value = 106545
$ws.Range("F33").Value = "Mindre märgborre"
$ws.Range("G33").Value = "Tomicus minor"
$ws.Range("H33").Value = "(Hartig, 1834)"
$ws.Range("M33").Value = "äldre gnagspår"
$ws.Range("Q33").Value = 485525
$ws.Range("R33").Value = 6783113
$ws.Range("S33").Value = 10

# Row 34
$ws.Range("A34").Value = 131202310
$ws.Range("Q34").Value = 485517
$ws.Range("R34").Value = 6783102
$ws.Range("AC34").ClearContents()

# Row 35
$ws.Range("B35").Value = 91830

# Row 38
$ws.Range("A38").Value = 131199044
$ws.Range("B38").Value = 57881
$ws.Range("E38").Value = 100049
$ws.Range("F38").Value = "Spillkråka"
$ws.Range("G38").Value = "Dryocopus martius"
$ws.Range("Q38").Value = 485494
$ws.Range("R38").Value = 6783163

# Row 39
$ws.Range("A39").Value = 131198860
$ws.Range("Q39").Value = 485442
$ws.Range("R39").Value = 6783096

# Row 40
$ws.Range("A40").Value = 131202630
$ws.Range("B40").Value = 57884
$ws.Range("E40").Value = 100109
$ws.Range("F40").Value = "Tretåig hackspett"
$ws.Range("G40").Value = "Picoides tridactylus"
$ws.Range("Q40").Value = 485526
$ws.Range("R40").Value = 6783143
